# Refresh this NATMI LR-pair sheet (Rspo3-Lrp6) with the new TPM-derived
# expression values for the "ECs" cluster. The ligand/receptor average &
# total expression values (columns G, H, M, N) were recomputed upstream
# from updated TPM data; the specificity and edge-weight columns
# (I, J, O, P, Q, R, S, T) are downstream derived metrics recomputed from
# those updated inputs across the whole sheet (rows 2-10).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.003420333333333333
$ws.Range("H2").Value = 0.010261
$ws.Range("I2").Value = 0.003549653112303053
$ws.Range("J2").Value = 0.003549653112303053
$ws.Range("M2").Value = 12.40685866666667
$ws.Range("N2").Value = 37.220576
$ws.Range("O2").Value = 0.1720325859617629
$ws.Range("P2").Value = 0.1720325859617629
$ws.Range("Q2").Value = 0.04243559225955555
$ws.Range("R2").Value = 0.381920330336
$ws.Range("S2").Value = 0.0006106560041767144
$ws.Range("T2").Value = 0.0006106560041767144

$ws.Range("G3").Value = 0.003420333333333333
$ws.Range("H3").Value = 0.010261
$ws.Range("I3").Value = 0.003549653112303053
$ws.Range("J3").Value = 0.003549653112303053
$ws.Range("O3").Value = 0.6097142007069145
$ws.Range("P3").Value = 0.6097142007069145
$ws.Range("Q3").Value = 0.1503993157541111
$ws.Range("R3").Value = 1.353593841787
$ws.Range("S3").Value = 0.002164273910154668
$ws.Range("T3").Value = 0.002164273910154668

$ws.Range("G4").Value = 0.003420333333333333
$ws.Range("H4").Value = 0.010261
$ws.Range("I4").Value = 0.003549653112303053
$ws.Range("J4").Value = 0.003549653112303053
$ws.Range("N4").Value = 47.220765
$ws.Range("O4").Value = 0.2182532133313226
$ws.Range("P4").Value = 0.2182532133313226
$ws.Range("Q4").Value = 0.05383691885166666
$ws.Range("R4").Value = 0.484532269665
$ws.Range("S4").Value = 0.0007747231979716716
$ws.Range("T4").Value = 0.0007747231979716716

$ws.Range("I5").Value = 0.3907064193682856
$ws.Range("J5").Value = 0.3907064193682855
$ws.Range("M5").Value = 12.40685866666667
$ws.Range("N5").Value = 37.220576
$ws.Range("O5").Value = 0.1720325859617629
$ws.Range("P5").Value = 0.1720325859617629
$ws.Range("Q5").Value = 4.670839031576889
$ws.Range("R5").Value = 42.03755128419201
$ws.Range("S5").Value = 0.06721423567578719
$ws.Range("T5").Value = 0.06721423567578719

$ws.Range("I6").Value = 0.3907064193682856
$ws.Range("J6").Value = 0.3907064193682855
$ws.Range("O6").Value = 0.6097142007069145
$ws.Range("P6").Value = 0.6097142007069145
$ws.Range("S6").Value = 0.2382192521961948
$ws.Range("T6").Value = 0.2382192521961947

$ws.Range("I7").Value = 0.3907064193682856
$ws.Range("J7").Value = 0.3907064193682855
$ws.Range("N7").Value = 47.220765
$ws.Range("O7").Value = 0.2182532133313226
$ws.Range("P7").Value = 0.2182532133313226
$ws.Range("Q7").Value = 5.925770527111667
$ws.Range("R7").Value = 53.331934744005
$ws.Range("S7").Value = 0.08527293149630363
$ws.Range("T7").Value = 0.08527293149630362

$ws.Range("I8").Value = 0.6057439275194114
$ws.Range("J8").Value = 0.6057439275194113
$ws.Range("M8").Value = 12.40685866666667
$ws.Range("N8").Value = 37.220576
$ws.Range("O8").Value = 0.1720325859617629
$ws.Range("P8").Value = 0.1720325859617629
$ws.Range("Q8").Value = 7.241581503505778
$ws.Range("R8").Value = 65.17423353155201
$ws.Range("S8").Value = 0.104207694281799
$ws.Range("T8").Value = 0.104207694281799

$ws.Range("I9").Value = 0.6057439275194114
$ws.Range("J9").Value = 0.6057439275194113
$ws.Range("O9").Value = 0.6097142007069145
$ws.Range("P9").Value = 0.6097142007069145
$ws.Range("S9").Value = 0.3693306746005651
$ws.Range("T9").Value = 0.369330674600565

$ws.Range("I10").Value = 0.6057439275194114
$ws.Range("J10").Value = 0.6057439275194113
$ws.Range("N10").Value = 47.220765
$ws.Range("O10").Value = 0.2182532133313226
$ws.Range("P10").Value = 0.2182532133313226
$ws.Range("Q10").Value = 9.187203830628334
$ws.Range("R10").Value = 82.68483447565501
$ws.Range("S10").Value = 0.1322055586370473
$ws.Range("T10").Value = 0.1322055586370473
